$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1695.75
$ws.Range("I19").Value = 873.8
$ws.Range("J19").Value = 2282.8572
$ws.Range("K19").Value = 873.8
$ws.Range("L19").Value = 2282.8572
$ws.Range("M19").Value = -698.8
$ws.Range("N19").Value = -2632.8572
$ws.Range("H135").Value = 164.75
$ws.Range("I135").Value = 86.333336
$ws.Range("K135").Value = 777.0000240000001
$ws.Range("M135").Value = 1757.999976
$ws.Range("H137").Value = 1281.6562
$ws.Range("I137").Value = 1115.7142
$ws.Range("J137").Value = 1598.4546
$ws.Range("K137").Value = 3347.1426
$ws.Range("L137").Value = 4795.3638
$ws.Range("M137").Value = -797.1425999999997
$ws.Range("N137").Value = -9895.363799999999
$ws.Range("H138").Value = 2044.64
$ws.Range("I138").Value = 956
$ws.Range("J138").Value = 2193.0908
$ws.Range("K138").Value = 2868
$ws.Range("L138").Value = 6579.2724
$ws.Range("M138").Value = 2272
$ws.Range("N138").Value = -16859.2724
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1976.0465
$ws.Range("I32").Value = 2111.8462
$ws.Range("K32").Value = 2111.8462
$ws.Range("M32").Value = -1824.8462
$ws.Range("H61").Value = 1497.1578
$ws.Range("I61").Value = 1307.375
$ws.Range("K61").Value = 1307.375
$ws.Range("M61").Value = -1095.375
$ws.Range("H74").Value = 597.9808
$ws.Range("I74").Value = 564.45
$ws.Range("J74").Value = 709.75
$ws.Range("K74").Value = 564.45
$ws.Range("L74").Value = 709.75
$ws.Range("M74").Value = 309.55
$ws.Range("N74").Value = -2457.75
$ws.Range("H77").Value = 597.9808
$ws.Range("I77").Value = 564.45
$ws.Range("J77").Value = 709.75
$ws.Range("K77").Value = 2822.25
$ws.Range("L77").Value = 3548.75
$ws.Range("M77").Value = 1545.75
$ws.Range("N77").Value = -12284.75
$ws.Range("H122").Value = 1565.8182
$ws.Range("I122").Value = 1551.2
$ws.Range("K122").Value = 4653.6
$ws.Range("M122").Value = -2203.6
$ws.Range("H136").Value = 1497.1578
$ws.Range("I136").Value = 1307.375
$ws.Range("K136").Value = 3922.125
$ws.Range("M136").Value = -1372.125
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 25001098
$ws.Range("I94").Value = 31250622
$ws.Range("K94").Value = 31250622
$ws.Range("M94").Value = -31250171
$ws.Range("H134").Value = 9524.654
$ws.Range("I134").Value = 6928.5264
$ws.Range("J134").Value = 16571.285
$ws.Range("K134").Value = 20785.5792
$ws.Range("L134").Value = 49713.855
$ws.Range("M134").Value = -18250.5792
$ws.Range("N134").Value = -54783.855
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 836.7846
$ws.Range("I31").Value = 758.87805
$ws.Range("J31").Value = 969.875
$ws.Range("K31").Value = 758.87805
$ws.Range("L31").Value = 969.875
$ws.Range("M31").Value = -463.87805
$ws.Range("N31").Value = -1559.875
$ws.Range("H34").Value = 836.7846
$ws.Range("I34").Value = 758.87805
$ws.Range("J34").Value = 969.875
$ws.Range("K34").Value = 758.87805
$ws.Range("L34").Value = 969.875
$ws.Range("M34").Value = -556.87805
$ws.Range("N34").Value = -1373.875
$ws.Range("H58").Value = 989.0952
$ws.Range("I58").Value = 1018.6923
$ws.Range("J58").Value = 941
$ws.Range("K58").Value = 1018.6923
$ws.Range("L58").Value = 941
$ws.Range("M58").Value = -815.6923
$ws.Range("N58").Value = -1347
$ws.Range("H99").Value = 1706.5217
$ws.Range("I99").Value = 1539.75
$ws.Range("J99").Value = 2087.7144
$ws.Range("K99").Value = 1539.75
$ws.Range("L99").Value = 2087.7144
$ws.Range("M99").Value = -41.75
$ws.Range("N99").Value = -5083.7144
$ws.Range("H126").Value = 1706.5217
$ws.Range("I126").Value = 1539.75
$ws.Range("J126").Value = 2087.7144
$ws.Range("K126").Value = 4619.25
$ws.Range("L126").Value = 6263.1432
$ws.Range("M126").Value = -2149.25
$ws.Range("N126").Value = -11203.1432
$ws.Range("H132").Value = 8872.875
$ws.Range("I132").Value = 18100.666
$ws.Range("K132").Value = 54301.99800000001
$ws.Range("M132").Value = -51771.99800000001
$ws.Range("H134").Value = 9260436
$ws.Range("I134").Value = 9804991
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 29414973
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -29412438
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 989.0952
$ws.Range("I136").Value = 1018.6923
$ws.Range("J136").Value = 941
$ws.Range("K136").Value = 3056.0769
$ws.Range("L136").Value = 2823
$ws.Range("M136").Value = -506.0769
$ws.Range("N136").Value = -7923
$ws.Range("H138").Value = 69780
$ws.Range("J138").Value = 69780
$ws.Range("L138").Value = 69780
$ws.Range("N138").Value = -80060
$ws.Range("H140").Value = 42450
$ws.Range("J140").Value = 42450
$ws.Range("L140").Value = 42450
$ws.Range("N140").Value = -52810
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1344.7333
$ws.Range("I68").Value = 660.2
$ws.Range("J68").Value = 2029.2667
$ws.Range("K68").Value = 1980.6
$ws.Range("L68").Value = 6087.800099999999
$ws.Range("M68").Value = -1169.6
$ws.Range("N68").Value = -7709.800099999999
$ws.Range("H71").Value = 1344.7333
$ws.Range("I71").Value = 660.2
$ws.Range("J71").Value = 2029.2667
$ws.Range("K71").Value = 5941.8
$ws.Range("L71").Value = 18263.4003
$ws.Range("M71").Value = -1885.8
$ws.Range("N71").Value = -26375.4003
$ws.Range("H76").Value = 6377.483
$ws.Range("I76").Value = 5182.6
$ws.Range("J76").Value = 6626.4165
$ws.Range("K76").Value = 15547.8
$ws.Range("L76").Value = 19879.2495
$ws.Range("M76").Value = -15164.8
$ws.Range("N76").Value = -20645.2495
$ws.Range("H79").Value = 6377.483
$ws.Range("I79").Value = 5182.6
$ws.Range("J79").Value = 6626.4165
$ws.Range("K79").Value = 15547.8
$ws.Range("L79").Value = 19879.2495
$ws.Range("M79").Value = -14221.8
$ws.Range("N79").Value = -22531.2495
$ws.Range("H99").Value = 2072
$ws.Range("J99").Value = 2771
$ws.Range("L99").Value = 8313
$ws.Range("N99").Value = -12805
$ws.Range("H107").Value = 4724.2
$ws.Range("J107").Value = 8558.308000000001
$ws.Range("L107").Value = 25674.924
$ws.Range("N107").Value = -29514.924
$ws.Range("H131").Value = 31251526
$ws.Range("J131").Value = 1867.0435
$ws.Range("L131").Value = 5601.1305
$ws.Range("N131").Value = -15681.1305
$ws.Range("H140").Value = 29440.775
$ws.Range("I140").Value = 43932.5
$ws.Range("J140").Value = 2527.5715
$ws.Range("K140").Value = 131797.5
$ws.Range("L140").Value = 7582.7145
$ws.Range("M140").Value = -126617.5
$ws.Range("N140").Value = -17942.7145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1534.4
$ws.Range("I102").Value = 1478.2667
$ws.Range("J102").Value = 1702.8
$ws.Range("K102").Value = 1478.2667
$ws.Range("L102").Value = 1702.8
$ws.Range("M102").Value = 143.7333000000001
$ws.Range("N102").Value = -4946.8
$ws.Range("H132").Value = 2173.4167
$ws.Range("I132").Value = 1766.6522
$ws.Range("K132").Value = 5299.9566
$ws.Range("M132").Value = -2769.9566
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2129
$ws.Range("J7").Value = 2149.25
$ws.Range("L7").Value = 2149.25
$ws.Range("N7").Value = -2373.25
$ws.Range("H40").Value = 2233
$ws.Range("I40").Value = 2237.6
$ws.Range("J40").Value = 2227.25
$ws.Range("K40").Value = 2237.6
$ws.Range("L40").Value = 2227.25
$ws.Range("M40").Value = -2101.6
$ws.Range("N40").Value = -2499.25
$ws.Range("H93").Value = 588.4286
$ws.Range("I93").Value = 548.5454999999999
$ws.Range("J93").Value = 734.6667
$ws.Range("K93").Value = 548.5454999999999
$ws.Range("L93").Value = 734.6667
$ws.Range("M93").Value = 699.4545000000001
$ws.Range("N93").Value = -3230.6667
$ws.Range("H122").Value = 15627115
$ws.Range("I122").Value = 27779686
$ws.Range("J122").Value = 2380
$ws.Range("K122").Value = 83339058
$ws.Range("L122").Value = 7140
$ws.Range("M122").Value = -83336608
$ws.Range("N122").Value = -12040
$ws.Range("H126").Value = 2129
$ws.Range("J126").Value = 2149.25
$ws.Range("L126").Value = 6447.75
$ws.Range("N126").Value = -11387.75
$ws.Range("H136").Value = 1607.4584
$ws.Range("I136").Value = 1647.5
$ws.Range("J136").Value = 1527.375
$ws.Range("K136").Value = 4942.5
$ws.Range("L136").Value = 4582.125
$ws.Range("M136").Value = -2392.5
$ws.Range("N136").Value = -9682.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1388.15
$ws.Range("I136").Value = 1221.1177
$ws.Range("K136").Value = 3663.3531
$ws.Range("M136").Value = -1113.3531
